# Apply the "May 9th" data update: the C:H sensor readings are shifted down by
# two rows relative to the timestamp column, two brand-new rows of sensor data
# are inserted at the top (using timestamps 0 and 100), and the timestamp
# series is extended all the way to row 31 (timestamp 2900), with the two
# rows that fall off the end of the old C:H series (old rows 20-21, now at
# timestamps 2000/2100) and eight further brand-new rows (2200-2900) getting
# freshly generated sensor data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 21   # last existing data row before the edit (timestamp 1900)

# Step 1: capture the existing C:H values (rows 2..21) before they get
# overwritten, so they can be re-written two rows further down (rows 4..23).
$oldCH = @{}
for ($r = 2; $r -le $lastDataRow; $r++) {
    $vals = @()
    for ($c = 3; $c -le 8; $c++) {
        $vals += $ws.Cells.Item($r, $c).Value2
    }
    $oldCH[$r] = $vals
}

# Step 2: re-write the old C:H values two rows further down (row r -> r+2).
for ($r = $lastDataRow; $r -ge 2; $r--) {
    $destRow = $r + 2
    $vals = $oldCH[$r]
    for ($i = 0; $i -lt 6; $i++) {
        $ws.Cells.Item($destRow, 3 + $i).Value2 = $vals[$i]
    }
}

# Step 3: new sensor-data rows inserted at the very top (rows 2 and 3),
# keeping timestamps 0 and 100.
$newTop = @(
    @(0,   0.2804546356201172, 0.4303635954856872, -0.691750168800354, 0.1050096067542932, -1.756468223065746, 0.4945203567645989),
    @(100, 0.1987819671630859, 0.2879692316055298, -0.9282988905906676, -0.0286234012063665, -0.7998002785809195, 0.0811297598541999)
)
for ($i = 0; $i -lt $newTop.Count; $i++) {
    $destRow = 2 + $i
    $vals = $newTop[$i]
    $ws.Cells.Item($destRow, 1).Value2 = $vals[0]
    $ws.Cells.Item($destRow, 2).Value = "struggle"
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($destRow, 3 + $j).Value2 = $vals[$j + 1]
    }
}

# Step 4: rows 22-31 (timestamps 2000-2900): two of these (2000, 2100) reuse
# the sensor data that used to belong to timestamps 1800/1900 (already moved
# into place by Step 2, at rows 22/23) and just need their timestamp/label
# cells filled in; the remaining eight (2200-2900) are entirely new rows.
$tailTimestamps = @(2000, 2100, 2200, 2300, 2400, 2500, 2600, 2700, 2800, 2900)
for ($i = 0; $i -lt $tailTimestamps.Count; $i++) {
    $destRow = 22 + $i
    $ws.Cells.Item($destRow, 1).Value2 = $tailTimestamps[$i]
    $ws.Cells.Item($destRow, 2).Value = "struggle"
}

$newBottom = @(
    @(-0.4514303207397461, -0.07753515243530271, -1.056098580360413, -0.4081483519807154, -0.6726997543354425, -0.2190668820118418),
    @(1.037992477416992, -1.273390769958496, 0.4362349510192871, 0.2211332225373814, 0.241335413285664, 0.08368853798934378),
    @(0.0754270553588867, 1.646718859672546, 1.695090532302856, 0.06768137718341787, 0.3379019900244107, 0.1505034766635118),
    @(-0.2560558319091797, 0.3026316165924072, -0.4233262538909912, 0.07254024853511698, 0.5556785336562575, -0.05807583201296457),
    @(0.6335611343383789, 0.8106564879417419, -1.443797469139099, 0.1816357883567719, 0.1322741392923868, -0.08515337003128903),
    @(0.09285736083984369, 0.7357764840126038, -1.646607518196106, -0.02734556931013958, -0.1169588795425942, 0.04497027853313797),
    @(0.0882749557495117, 0.1726978719234466, -0.9354652166366576, -0.02540700723017953, -0.06986615411481072, -0.074921377335808),
    @(0.2656211853027344, 0.4902379512786865, -0.8409426212310791, 0.02237761537639455, -0.07008743807863513, -0.003453258577050004)
)
for ($i = 0; $i -lt $newBottom.Count; $i++) {
    $destRow = 24 + $i
    $vals = $newBottom[$i]
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($destRow, 3 + $j).Value2 = $vals[$j]
    }
}
